$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Prepare formatting for the two brand-new cells that will appear
#    in column L (L6, L7) once columns M:Q are removed. They need the
#    same "positive number" style (fillId green) that is already used
#    by C3 (style index 7 in the source workbook).
# ------------------------------------------------------------------
$ws.Range("C3").Copy()
$ws.Range("L6").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("L7").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 2) Drop the last 5 years (columns M:Q) - the model now only spans
#    10 years (columns B:L) instead of 15 (columns B:Q).
# ------------------------------------------------------------------
$ws.Range("M1:Q15").Clear()

# ------------------------------------------------------------------
# 3) Depreciation (row 3): 20000 -> 30000 for years 1-10 (C:L)
# ------------------------------------------------------------------
$ws.Range("C3:L3").Value = 30000

# ------------------------------------------------------------------
# 4) Incoming Payments (row 4): 1050000 -> 700000 for years 1-10 (C:L)
# ------------------------------------------------------------------
$ws.Range("C4:L4").Value = 700000

# ------------------------------------------------------------------
# 5) Outgoing Payments (row 5): year-0 value doubles
# ------------------------------------------------------------------
$ws.Range("B5").Value = -140000

# ------------------------------------------------------------------
# 6) Residual (row 6): now realized in year 10 (L) instead of 15 (Q)
# ------------------------------------------------------------------
$ws.Range("L6").Value = 140000

# ------------------------------------------------------------------
# 7) restricted Equity (row 7): value doubles, now realized in L
# ------------------------------------------------------------------
$ws.Range("B7").Value = -300000
$ws.Range("L7").Value = 300000

# ------------------------------------------------------------------
# 8) Yearly Net (row 8)
# ------------------------------------------------------------------
$ws.Range("B8").Value = -1440000
$ws.Range("C8:K8").Value = 380000
$ws.Range("L8").Value = 820000

# ------------------------------------------------------------------
# 9) Present Value (row 9)
# ------------------------------------------------------------------
$ws.Range("B9").Value = -1440000
$ws.Range("C9").Value = 350553.5055350553
$ws.Range("D9").Value = 323388.842744516
$ws.Range("E9").Value = 298329.1907237232
$ws.Range("F9").Value = 275211.4305569402
$ws.Range("G9").Value = 253885.0835396127
$ws.Range("H9").Value = 234211.3316786095
$ws.Range("I9").Value = 216062.1140946582
$ws.Range("J9").Value = 199319.2934452566
$ws.Range("K9").Value = 183873.8869421186
$ws.Range("L9").Value = 366033.6650139281

# ------------------------------------------------------------------
# 10) Accumulated Present Value (row 10). D10:F10 now turn negative so
#     they must switch from the "positive" style to the "negative" one
#     (copy the negative-number format already present on C10).
# ------------------------------------------------------------------
$ws.Range("C10").Copy()
$ws.Range("D10:F10").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B10").Value = -1440000
$ws.Range("C10").Value = -1089446.494464945
$ws.Range("D10").Value = -766057.6517204286
$ws.Range("E10").Value = -467728.4609967054
$ws.Range("F10").Value = -192517.0304397652
$ws.Range("G10").Value = 61368.05309984752
$ws.Range("H10").Value = 295579.384778457
$ws.Range("I10").Value = 511641.4988731152
$ws.Range("J10").Value = 710960.7923183718
$ws.Range("K10").Value = 894834.6792604905
$ws.Range("L10").Value = 1260868.344274419

# ------------------------------------------------------------------
# 11) Net Present Value (row 11)
# ------------------------------------------------------------------
$ws.Range("B11").Value = 1260868.344274419
